$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A3 was stored as text "44189151" — turn it into a real number, matching
# the numeric DNI cells elsewhere in the column.
$ws.Range("A3").Value = 44189151

# Helper pattern used below for the Fecha/Hora (and row-6 DNI) cells: force
# Text number format so the literal string isn't reinterpreted as a date or
# serial number, assign the value, then clear the format again so the cell
# doesn't end up carrying an explicit style index.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# New row 4 — first scan for Bazan, Ruth
$ws.Range("A4").Value = 37128854
$ws.Range("B4").Value = "Bazan"
$ws.Range("C4").Value = "Ruth"
Set-TextValue $ws.Range("D4") "2023-11-13"
Set-TextValue $ws.Range("E4") "08:01:57"

# New row 5 — second scan for Bazan, Ruth
$ws.Range("A5").Value = 37128854
$ws.Range("B5").Value = "Bazan"
$ws.Range("C5").Value = "Ruth"
Set-TextValue $ws.Range("D5") "2023-11-13"
Set-TextValue $ws.Range("E5") "08:06:10"

# New row 6 — Guiñazu, Alejandro (DNI stays text here, unlike rows 4/5)
Set-TextValue $ws.Range("A6") "27775770"
$ws.Range("B6").Value = "Guiñazu"
$ws.Range("C6").Value = "Alejandro"
Set-TextValue $ws.Range("D6") "2023-11-13"
Set-TextValue $ws.Range("E6") "08:08:40"
